# "Added external pulls and decaps"
#
# TPcape_BOM sheet: populate VALUE/DESCRIPTION for two component arrays that
# previously had no part info filled in:
#   - row 37: RA17, RA21, RA23a, RA24 external pull-up resistor array
#   - row 25: CA1, CA2, CA3 decoupling capacitor array
#
# Connectors sheet: a new (blank) column was inserted before column G.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TPcape_BOM")
$ws2 = $wb.Worksheets.Item("Connectors")
$ws3 = $wb.Worksheets.Item("interface")

# --- TPcape_BOM: fill in the new component rows ---------------------------

# E19 already carries the "description" text style (Arial 9pt) used
# elsewhere in this column; copy that formatting onto the two new
# description cells before filling them in.
$ws1.Range("E19").Copy()
$ws1.Range("E37").PasteSpecial(-4122)
$ws1.Range("E25").PasteSpecial(-4122)

# Row 37 - external pull-up resistor array (RA17, RA21, RA23a, RA24)
$ws1.Range("E37").Value = "RES ARRAY 22K OHM 4 RES 1206"
$ws1.Range("D37").Value = "22K"

# Row 25 - decoupling capacitor array (CA1, CA2, CA3)
$ws1.Range("D25").Value = "1000pF"
$ws1.Range("E25").Value = "CAP ARRAY 1000PF 100V X7R 1206"

# --- Connectors: insert a blank column before column G --------------------
$ws2.Range("G1").EntireColumn.Insert()

# --- Restore the editing session's cursor / active-sheet state ------------
$ws1.Activate() | Out-Null
$ws1.Range("E28").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A55").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F18").Select() | Out-Null
